$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "Single Match Simulator" sheet: add a new "Is Singles?" setting row (7)
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Single Match Simulator")

# Row 7 was previously blank (gap between row 6 and row 8). Fill it in and
# copy the input-cell style from B6 onto B7 so it matches the other setting
# cells (fill + border).
$ws1.Range("A7").Value = "Is Singles?"

$ws1.Range("B6").Copy()
$ws1.Range("B7").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("B7").Value = "No"

$ws1.Range("C7").Value = "1v1 matches count 0.5x weight"

# New dropdown validation on B7, matching B6's Yes/No list validation.
$ws1.Range("B7").Validation.Add(3, 1, 1, """Yes,No""")
$ws1.Range("B7").Validation.IgnoreBlank = $true
$ws1.Range("B7").Validation.InCellDropdown = $true
$ws1.Range("B7").Validation.ShowInput = $false
$ws1.Range("B7").Validation.ShowError = $false

# Update the Match Weight formula (L14) to also factor in the new Singles
# setting: 0.5x when B7 = "Yes", otherwise unchanged behavior.
$ws1.Range("L14").Formula = '=IF(B7="Yes", 0.5, 1) * IF(B6="Yes", 1, IF(B4="Sets", IF(MAX(B8,B9)>=6, 1, 0.5), IF(B5>21, 1, 0.5)))'

# ---------------------------------------------------------------------------
# 2) "Documentation" sheet: insert a new "Singles Match Weight" row before
#    the existing "Player Weight" row, shifting everything below down by one.
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Documentation")

$ws2.Rows.Item(7).Insert()

# Copy formatting from the row below (now the old "Player Weight" row, which
# has already shifted to row 8) so the new row matches the surrounding style.
$ws2.Range("A8:B8").Copy()
$ws2.Range("A7:B7").PasteSpecial(-4122)  # xlPasteFormats

$ws2.Range("A7").Value = "Singles Match Weight"
$ws2.Range("B7").Value = "1v1 matches count for 0.5x weight."

# Update the Delta (ELO) formula documentation (now row 21) to include the
# new SinglesWeight factor.
$ws2.Range("B21").Formula = "=ROUND(K * MarginMult * MatchWeight * SinglesWeight * EffWeight * (Result - Expected), 0)"
